$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")
$ws.Activate()

# Insert a new (blank) column in front of column N. This pushes the old
# N/O/P columns ("Late" header + its data, and "Outstanding" header + its
# data) one column to the right, becoming O/P/Q - matching the new
# "Variable Instalments" layout used for RBI loans.
$ws.Columns("N:N").Insert()

# Match the column width Excel stores for the freshly inserted column
# (same stored width, 10 characters, as column M).
$ws.Columns("N:N").ColumnWidth = 9.17

# Update the active cell/selection to match the saved state.
$ws.Range("U7").Select() | Out-Null
